{"js": "const pairs = [\n  [\"2025-07-12 Saturday\", \"2025-07-13 Sunday\"],\n  [\"31\u00d791=2821\", \"21\u00d724=504\"],\n  [\"42\u00d713=546\", \"43\u00d769=2967\"],\n  [\"59\u00d753=3127\", \"47\u00d794=4418\"],\n  [\"88\u00d787=7656\", \"19\u00d722=418\"],\n  [\"97\u00d768=6596\", \"89\u00d746=4094\"],\n  [\"54\u00d740=2160\", \"71\u00d754=3834\"],\n  [\"27\u00d729=783\", \"88\u00d774=6512\"],\n  [\"53\u00d764=3392\", \"46\u00d753=2438\"],\n  [\"98\u00d780=7840\", \"83\u00d799=8217\"],\n  [\"27\u00d759=1593\", \"76\u00d795=7220\"],\n  [\"44\u00d792=4048\", \"37\u00d743=1591\"],\n  [\"79\u00d716=1264\", \"58\u00d752=3016\"],\n  [\"13\u00d723=299\", \"44\u00d794=4136\"],\n  [\"79\u00d787=6873\", \"15\u00d749=735\"],\n  [\"96\u00d753=5088\", \"88\u00d783=7304\"],\n  [\"95\u00d783=7885\", \"78\u00d733=2574\"],\n  [\"20\u00d717=340\", \"70\u00d764=4480\"],\n  [\"54\u00d739=2106\", \"96\u00d787=8352\"],\n  [\"20\u00d718=360\", \"54\u00d732=1728\"],\n  [\"20\u00d744=880\", \"68\u00d731=2108\"],\n  [\"35\u00d746=1610\", \"65\u00d759=3835\"],\n  [\"14\u00d727=378\", \"41\u00d747=1927\"],\n  [\"90\u00d713=1170\", \"94\u00d771=6674\"],\n  [\"89\u00d759=5251\", \"30\u00d719=570\"],\n  [\"55\u00d754=2970\", \"41\u00d753=2173\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$pairs = @(\n  @(\"2025-07-12 Saturday\", \"2025-07-13 Sunday\"),\n  @(\"31\u00d791=2821\", \"21\u00d724=504\"),\n  @(\"42\u00d713=546\", \"43\u00d769=2967\"),\n  @(\"59\u00d753=3127\", \"47\u00d794=4418\"),\n  @(\"88\u00d787=7656\", \"19\u00d722=418\"),\n  @(\"97\u00d768=6596\", \"89\u00d746=4094\"),\n  @(\"54\u00d740=2160\", \"71\u00d754=3834\"),\n  @(\"27\u00d729=783\", \"88\u00d774=6512\"),\n  @(\"53\u00d764=3392\", \"46\u00d753=2438\"),\n  @(\"98\u00d780=7840\", \"83\u00d799=8217\"),\n  @(\"27\u00d759=1593\", \"76\u00d795=7220\"),\n  @(\"44\u00d792=4048\", \"37\u00d743=1591\"),\n  @(\"79\u00d716=1264\", \"58\u00d752=3016\"),\n  @(\"13\u00d723=299\", \"44\u00d794=4136\"),\n  @(\"79\u00d787=6873\", \"15\u00d749=735\"),\n  @(\"96\u00d753=5088\", \"88\u00d783=7304\"),\n  @(\"95\u00d783=7885\", \"78\u00d733=2574\"),\n  @(\"20\u00d717=340\", \"70\u00d764=4480\"),\n  @(\"54\u00d739=2106\", \"96\u00d787=8352\"),\n  @(\"20\u00d718=360\", \"54\u00d732=1728\"),\n  @(\"20\u00d744=880\", \"68\u00d731=2108\"),\n  @(\"35\u00d746=1610\", \"65\u00d759=3835\"),\n  @(\"14\u00d727=378\", \"41\u00d747=1927\"),\n  @(\"90\u00d713=1170\", \"94\u00d771=6674\"),\n  @(\"89\u00d759=5251\", \"30\u00d719=570\"),\n  @(\"55\u00d754=2970\", \"41\u00d753=2173\"),\n)\n\n$d = $word.ActiveDocument\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $p[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $p[1]\n  $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}"}
